$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells - copy formatting from an existing header cell, then set values
$ws.Range("A1").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122)
$ws.Range("G1").Value = "Elapsed Time"
$ws.Range("H1").Value = "CPU"

# Update existing MSE/R2/MAE values
$ws.Range("B2").Value = 0.05251846045315886
$ws.Range("C2").Value = 0.9984555216771278
$ws.Range("D2").Value = 0.1681329394362654

$ws.Range("B3").Value = 0.06090733564894349
$ws.Range("C3").Value = 0.999422266625869
$ws.Range("D3").Value = 0.1782011012841147

$ws.Range("B4").Value = 0.05576529082056546
$ws.Range("C4").Value = 0.9992449723851116
$ws.Range("D4").Value = 0.1916073009001543

# New Elapsed Time / CPU columns
$ws.Range("G2").Value = 0.3776785511166963
$ws.Range("H2").Value = 0.968

$ws.Range("G3").Value = 0.3776785511166963
$ws.Range("H3").Value = 0.968

$ws.Range("G4").Value = 0.3776785511166963
$ws.Range("H4").Value = 0.968
